$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number for every data row (2-484).
# The value changes from 45186 to 45188 for all of these rows.
$ws.Range("C2:C484").Value = 45188
